# Generate Report for Handoff
#
# The localization-status report tracks 3 source files:
#   569f2bbe-ada9-402d-892b-0544b541c87f.md
#   ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md
#   ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md
#
# 569f2bbe... has just been hooded off again ("Ready for handoff"), so it is
# regenerated (moved) to the bottom of each table, while the other two rows
# shift up one position (their data otherwise unchanged). This script rebuilds
# the hyperlinks + cell values on all three sheets to reflect the new report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", "", "", "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md", "", "", "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/.localization-config", "", "", ".localization-config")

$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", "", "", "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d9f9969a56b01b4aeabcf578aea4ba828c712d0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf")
$ws.Range("D2").Value = "2016-03-09 14:29:08"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/33863e781d44476cfabc240ac916a1e93a69defe/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bdf0a7116025d70c266d8f38770b9faef2c21759/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf")
$ws.Range("G2").Value = "2016-03-09 14:29:36"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("H2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md", "", "", "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8936cbf64460a8d88a1cbecad0156059d56b583f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf")
$ws.Range("D3").Value = "2016-03-09 14:29:08"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bfb945533f82a05b9749d4a053efda92fea8d709/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/344e3e673f5f9c98f8f96d94a512be8902b61e14/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf")
$ws.Range("G3").Value = "2016-03-09 14:29:36"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("H3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8936cbf64460a8d88a1cbecad0156059d56b583f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf")
$ws.Range("D4").Value = "2016-03-09 14:31:25"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bfb945533f82a05b9749d4a053efda92fea8d709/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/344e3e673f5f9c98f8f96d94a512be8902b61e14/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf")
$ws.Range("G4").Value = "2016-03-09 14:30:39"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("H4").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/.localization-config", "", "", ".localization-config")
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", "", "", "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce7e6726a0ed024ca3bc9b159c8c52e91cebb596/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf")
$ws.Range("D2").Value = "2016-03-09 14:29:11"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6c70c04c6bc032894d0d3bca851bfc4f941627e2/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/56b34ff25c388febf7e448993b4d8361cf9e29af/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf")
$ws.Range("G2").Value = "2016-03-09 14:29:41"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("H2").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md", "", "", "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9294d13ee7f6855343eb677d18ff7b6b8dc09f84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf")
$ws.Range("D3").Value = "2016-03-09 14:29:11"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a932d76745e8ce1a0dc1da24c37fe5b6aa9d96c8/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b206085512fda197301b56a2f586d7b0b06e0cb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf")
$ws.Range("G3").Value = "2016-03-09 14:29:41"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("H3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/e2e/ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9294d13ee7f6855343eb677d18ff7b6b8dc09f84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf")
$ws.Range("D4").Value = "2016-03-09 14:31:28"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a932d76745e8ce1a0dc1da24c37fe5b6aa9d96c8/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b206085512fda197301b56a2f586d7b0b06e0cb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", "", "", "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf")
$ws.Range("G4").Value = "2016-03-09 14:30:50"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("H4").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9b1f75d111502b7de4546a88223cd79bbaeb4103/.localization-config", "", "", ".localization-config")
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"
